# Weekly update: insert two new price-report rows for "Betarraga" at the
# top of the data block (rows 1319-1320), pushing the existing 1319-1382
# rows down to 1321-1384.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 1319 (existing data shifts
# down to 1321:1384; Excel copies the row-above formatting, matching the
# style="2" date format already used by column D).
$ws.Range("A1319:A1320").EntireRow.Insert()

# New row 1319 - "Primera" category
$ws.Cells.Item(1319, 1).Value = 6
$ws.Cells.Item(1319, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1319, 3).Value = "Metropolitana"
$ws.Cells.Item(1319, 4).Value = 45041
$ws.Cells.Item(1319, 5).Value = 13
$ws.Cells.Item(1319, 6).Value = 100114014
$ws.Cells.Item(1319, 7).Value = "Betarraga"
$ws.Cells.Item(1319, 8).Value = "Sin especificar"
$ws.Cells.Item(1319, 9).Value = "Primera"
$ws.Cells.Item(1319, 10).Value = 33000
$ws.Cells.Item(1319, 11).Value = 90
$ws.Cells.Item(1319, 12).Value = 90
$ws.Cells.Item(1319, 13).Value = 90
$ws.Cells.Item(1319, 14).Value = "`$/unidad"
$ws.Cells.Item(1319, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1319, 16).Value = 90
$ws.Cells.Item(1319, 17).Value = 1
$ws.Cells.Item(1319, 18).Value = "Hortaliza"

# New row 1320 - "Segunda" category
$ws.Cells.Item(1320, 1).Value = 6
$ws.Cells.Item(1320, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1320, 3).Value = "Metropolitana"
$ws.Cells.Item(1320, 4).Value = 45041
$ws.Cells.Item(1320, 5).Value = 13
$ws.Cells.Item(1320, 6).Value = 100114014
$ws.Cells.Item(1320, 7).Value = "Betarraga"
$ws.Cells.Item(1320, 8).Value = "Sin especificar"
$ws.Cells.Item(1320, 9).Value = "Segunda"
$ws.Cells.Item(1320, 10).Value = 22000
$ws.Cells.Item(1320, 11).Value = 65
$ws.Cells.Item(1320, 12).Value = 65
$ws.Cells.Item(1320, 13).Value = 65
$ws.Cells.Item(1320, 14).Value = "`$/unidad"
$ws.Cells.Item(1320, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1320, 16).Value = 65
$ws.Cells.Item(1320, 17).Value = 1
$ws.Cells.Item(1320, 18).Value = "Hortaliza"
